$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '71.471.54'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '3.875.95'
$ws.Range('E3').Value = '  -2.65%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '601.88'
$ws.Range('E5').Value = '  -1.64%  '
Set-TextValue $ws.Range('D6') '172.79'
$ws.Range('E6').Value = '  +5.65%  '
Set-TextValue $ws.Range('D7') '0.669'
$ws.Range('E7').Value = '  -2.13%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -0.79%  '
Set-TextValue $ws.Range('D10') '0.177'
$ws.Range('E10').Value = '  +5.41%  '
Set-TextValue $ws.Range('D11') '53.93'
$ws.Range('E11').Value = '  -0.87%  '
Set-TextValue $ws.Range('D12') '0.0000323'
$ws.Range('E12').Value = '  +0.80%  '
Set-TextValue $ws.Range('D13') '11.54'
$ws.Range('E13').Value = '  +5.09%  '
$ws.Range('D14').Value = '4.493.48'
$ws.Range('E14').Value = '  -2.67%  '
Set-TextValue $ws.Range('D15') '21.20'
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').Value = '3.877.49'
$ws.Range('E16').Value = '  -2.92%  '
Set-TextValue $ws.Range('D17') '13.96'
$ws.Range('E17').Value = '  -1.39%  '
$ws.Range('E18').Value = '  -3.74%  '
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').Value = '71.286.31'
$ws.Range('E20').Value = '  -1.76%  '
Set-TextValue $ws.Range('D21') '440.61'
$ws.Range('E21').Value = '  -0.03%  '
Set-TextValue $ws.Range('D22') '4.79'
$ws.Range('E22').Value = '  -1.69%  '
Set-TextValue $ws.Range('D23') '94.67'
$ws.Range('E23').Value = '  -2.09%  '
$ws.Range('E24').Value = '  -4.55%  '
Set-TextValue $ws.Range('D25') '13.90'
$ws.Range('E25').Value = '  -3.68%  '
Set-TextValue $ws.Range('D26') '11.79'
$ws.Range('E26').Value = '  +2.59%  '
$ws.Range('E27').Value = '  -6.63%  '
$ws.Range('E28').Value = '  +0.21%  '
Set-TextValue $ws.Range('D29') '10.49'
$ws.Range('E29').Value = '  -0.55%  '
Set-TextValue $ws.Range('D30') '8.67'
$ws.Range('E30').Value = '  +11.15%  '
Set-TextValue $ws.Range('D31') '35.26'
$ws.Range('E31').Value = '  -3.43%  '
$ws.Range('E32').Value = '  -2.56%  '
Set-TextValue $ws.Range('D33') '47.97'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('E34').Value = '  -3.72%  '
Set-TextValue $ws.Range('D35') '0.0000100'
$ws.Range('E35').Value = '  +11.19%  '
Set-TextValue $ws.Range('D36') '69.63'
$ws.Range('E36').Value = '  -2.57%  '
Set-TextValue $ws.Range('D37') '631.87'
$ws.Range('E37').Value = '  -4.16%  '
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('E39').Value = '  +0.45%  '
Set-TextValue $ws.Range('D40') '0.999'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('E41').Value = '  -0.18%  '
Set-TextValue $ws.Range('D42') '3.29'
$ws.Range('E42').Value = '  -2.45%  '
Set-TextValue $ws.Range('D43') '2.88'
$ws.Range('E43').Value = '  +8.62%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D44') '3.18'
$ws.Range('E44').Value = '  +19.40%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D45') '0.0473'
$ws.Range('E45').Value = '  -3.55%  '
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('E47').Value = '  -12.40%  '
$ws.Range('E48').Value = '  -3.84%  '
$ws.Range('D49').Value = '2.921.69'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('E51').Value = '  +2.94%  '
